$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.537.60'
$ws.Range('E2').Value = '  +6.32%  '
$ws.Range('D3').Value = '3.510.33'
$ws.Range('E3').Value = '  +10.24%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '561.75'
$ws.Range('E5').Value = '  +10.41%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '185.67'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +10.90%  '
$ws.Range('E7').Value = '  +10.49%  '
$ws.Range('D8').Value = '3.505.69'
$ws.Range('E8').Value = '  +10.17%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('E10').Value = '  +10.11%  '
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  +21.85%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '55.47'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +8.70%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000280'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +13.00%  '
$ws.Range('E14').Value = '  +8.42%  '
$ws.Range('D15').Value = '4.075.77'
$ws.Range('E15').Value = '  +9.18%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '18.76'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +11.82%  '
$ws.Range('D17').Value = '3.514.06'
$ws.Range('E17').Value = '  +9.27%  '
$ws.Range('E18').Value = '  +6.97%  '
$ws.Range('D19').Value = '66.653.41'
$ws.Range('E19').Value = '  +7.10%  '
$ws.Range('D20').Value = '12.08'
$ws.Range('E20').Value = '  +12.27%  '
$ws.Range('E21').Value = '  +8.78%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '419.82'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +15.42%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.10'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +16.09%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '85.76'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +8.68%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.16'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.03'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.59%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.91'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +12.74%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '12.36'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +14.62%  '
$ws.Range('E29').Value = '  +0.67%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.16'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +17.21%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '30.44'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +10.19%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '6.70'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +5.29%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '623.79'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('E34').Value = '  +9.43%  '
$ws.Range('E35').Value = '  +10.80%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '59.97'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +6.68%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.150'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +25.43%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0824'
$ws.Range('E38').Value = '  +18.96%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '38.22'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +11.06%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('E41').Value = '  +7.08%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.150.01'
$ws.Range('E42').Value = '  +13.01%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.36'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +15.06%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.25%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.65'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +5.20%  '
$ws.Range('E46').Value = '  +15.68%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0419'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +10.12%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.28'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +12.28%  '
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  +10.21%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '139.82'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +4.12%  '
